# Insert two new weekly rows of data right before the current row 304
# (new "Region Metropolitana" readings dated 44798), pushing the existing
# rows 304-363 down to 306-365.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 304 - everything that was row 304
# onward shifts down by two rows (to 306 onward).
$ws.Rows(304).Insert()
$ws.Rows(304).Insert()

# --- New row 304 ("Primera" quality) ---
$ws.Range("A304").Value2 = 11
$ws.Range("B304").Value2 = 'Vega Monumental Concepción'
$ws.Range("C304").Value2 = 'Bíobío'
$ws.Range("D304").Value2 = 44798
$ws.Range("E304").Value2 = 8
$ws.Range("F304").Value2 = 100114014
$ws.Range("G304").Value2 = 'Betarraga'
$ws.Range("H304").Value2 = 'Sin especificar'
$ws.Range("I304").Value2 = 'Primera'
$ws.Range("J304").Value2 = 600
$ws.Range("K304").Value2 = 800
$ws.Range("L304").Value2 = 900
$ws.Range("M304").Value2 = 850
$ws.Range("N304").Value2 = '$/paquete 5 unidades'
$ws.Range("O304").Value2 = 'Región Metropolitana'
$ws.Range("P304").Value2 = 170
$ws.Range("Q304").Value2 = 5
$ws.Range("R304").Value2 = 'Hortaliza'

# --- New row 305 ("Segunda" quality) ---
$ws.Range("A305").Value2 = 11
$ws.Range("B305").Value2 = 'Vega Monumental Concepción'
$ws.Range("C305").Value2 = 'Bíobío'
$ws.Range("D305").Value2 = 44798
$ws.Range("E305").Value2 = 8
$ws.Range("F305").Value2 = 100114014
$ws.Range("G305").Value2 = 'Betarraga'
$ws.Range("H305").Value2 = 'Sin especificar'
$ws.Range("I305").Value2 = 'Segunda'
$ws.Range("J305").Value2 = 300
$ws.Range("K305").Value2 = 700
$ws.Range("L305").Value2 = 700
$ws.Range("M305").Value2 = 700
$ws.Range("N305").Value2 = '$/paquete 5 unidades'
$ws.Range("O305").Value2 = 'Región Metropolitana'
$ws.Range("P305").Value2 = 140
$ws.Range("Q305").Value2 = 5
$ws.Range("R305").Value2 = 'Hortaliza'
